# Apply updated odds values to Sheet1, as described in the commit
# "Atualizando o arquivo XLSX" (updating the XLSX file).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("AA5").Value = 6.2
$ws.Range("AB5").Value = 15.5
$ws.Range("AC5").Value = 80
$ws.Range("AD5").Value = 7.4
$ws.Range("AG5").Value = 28
$ws.Range("AI5").Value = 35
$ws.Range("G5").Value = 3
$ws.Range("I5").Value = 2.42
$ws.Range("P5").Value = 1.45
$ws.Range("Q5").Value = 2.62
$ws.Range("R5").Value = 1.82
$ws.Range("S5").Value = 1.88
$ws.Range("V5").Value = 11.5
$ws.Range("X5").Value = 30
$ws.Range("Y5").Value = 45
# Row 6
$ws.Range("AE6").Value = 15
$ws.Range("G6").Value = 2.5
$ws.Range("H6").Value = 2.63
$ws.Range("I6").Value = 3.4
$ws.Range("J6").Value = 1.17
$ws.Range("K6").Value = 5
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = 10
$ws.Range("V6").Value = 11
$ws.Range("W6").Value = 26
$ws.Range("X6").Value = 26
# Row 7
$ws.Range("AA7").Value = 6.1
$ws.Range("AB7").Value = 19
$ws.Range("AC7").Value = 120
$ws.Range("AD7").Value = 6.1
$ws.Range("AE7").Value = 9.75
$ws.Range("AF7").Value = 10
$ws.Range("AG7").Value = 23
$ws.Range("AH7").Value = 23
$ws.Range("AI7").Value = 45
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 2.32
$ws.Range("L7").Value = 1.47
$ws.Range("M7").Value = 2.35
$ws.Range("N7").Value = 2.35
$ws.Range("O7").Value = 1.47
$ws.Range("P7").Value = 1.53
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 2.05
$ws.Range("S7").Value = 1.62
$ws.Range("T7").Value = 7.1
$ws.Range("V7").Value = 11.75
$ws.Range("W7").Value = 37
$ws.Range("X7").Value = 32
$ws.Range("Y7").Value = 50
$ws.Range("Z7").Value = 6.9
# Row 8
$ws.Range("AB8").Value = 14.5
$ws.Range("AC8").Value = 70
$ws.Range("AD8").Value = 10.25
$ws.Range("AE8").Value = 19.5
$ws.Range("AF8").Value = 12.5
$ws.Range("AH8").Value = 35
$ws.Range("AI8").Value = 40
$ws.Range("L8").Value = 1.3
$ws.Range("M8").Value = 2.92
$ws.Range("N8").Value = 1.88
$ws.Range("O8").Value = 1.72
$ws.Range("P8").Value = 1.42
$ws.Range("Q8").Value = 2.47
$ws.Range("R8").Value = 1.72
$ws.Range("S8").Value = 1.88
$ws.Range("T8").Value = 7.2
$ws.Range("V8").Value = 8.25
$ws.Range("W8").Value = 17.5
$ws.Range("X8").Value = 15.5
$ws.Range("Z8").Value = 9.25
# Row 9
$ws.Range("AD9").Value = 7.9
$ws.Range("AE9").Value = 14.5
$ws.Range("AF9").Value = 11
$ws.Range("AH9").Value = 29
$ws.Range("AI9").Value = 40
$ws.Range("AJ9").Value = 800
$ws.Range("H9").Value = 2.95
$ws.Range("I9").Value = 3
$ws.Range("T9").Value = 6.5
$ws.Range("U9").Value = 10.5
$ws.Range("X9").Value = 23
# Row 11
$ws.Range("AA11").Value = 7.9
$ws.Range("AD11").Value = 16
$ws.Range("AE11").Value = 29
$ws.Range("H11").Value = 3.95
$ws.Range("I11").Value = 4.5
$ws.Range("M11").Value = 3.75
$ws.Range("T11").Value = 8.5
$ws.Range("Z11").Value = 14
# Row 17
$ws.Range("AA17").Value = 7
$ws.Range("AD17").Value = 11
$ws.Range("AF17").Value = 17
$ws.Range("G17").Value = 1.65
$ws.Range("I17").Value = 5
$ws.Range("K17").Value = 8.5
$ws.Range("N17").Value = 2.15
$ws.Range("O17").Value = 1.67
$ws.Range("T17").Value = 6
# Row 19
$ws.Range("AD19").Value = 9.5
$ws.Range("AH19").Value = 17
$ws.Range("K19").Value = 13
$ws.Range("L19").Value = 1.22
$ws.Range("M19").Value = 4
$ws.Range("N19").Value = 1.75
$ws.Range("O19").Value = 2.05
$ws.Range("T19").Value = 11
# Row 20
$ws.Range("AD20").Value = 9.5
$ws.Range("H20").Value = 3.25
$ws.Range("L20").Value = 1.4
$ws.Range("M20").Value = 2.75
$ws.Range("P20").Value = 1.47
$ws.Range("V20").Value = 9
$ws.Range("Z20").Value = 7.5
# Row 22
$ws.Range("AA22").Value = 6.6
$ws.Range("AB22").Value = 16.5
$ws.Range("AC22").Value = 90
$ws.Range("AD22").Value = 10.25
$ws.Range("AE22").Value = 21
$ws.Range("AF22").Value = 13.5
$ws.Range("AG22").Value = 65
$ws.Range("AH22").Value = 40
$ws.Range("AI22").Value = 50
$ws.Range("AJ22").Value = 800
$ws.Range("G22").Value = 1.85
$ws.Range("H22").Value = 3.4
$ws.Range("I22").Value = 4.05
$ws.Range("J22").Value = 1.07
$ws.Range("K22").Value = 6.8
$ws.Range("L22").Value = 1.35
$ws.Range("M22").Value = 2.95
$ws.Range("N22").Value = 2.02
$ws.Range("O22").Value = 1.7
$ws.Range("P22").Value = 1.45
$ws.Range("Q22").Value = 2.55
$ws.Range("R22").Value = 1.9
$ws.Range("S22").Value = 1.8
$ws.Range("T22").Value = 6.5
$ws.Range("U22").Value = 8.25
$ws.Range("V22").Value = 8.5
$ws.Range("W22").Value = 15
$ws.Range("X22").Value = 15.5
$ws.Range("Y22").Value = 30
$ws.Range("Z22").Value = 6.8
# Row 29
$ws.Range("AA29").Value = 6.1
$ws.Range("AB29").Value = 15
$ws.Range("AC29").Value = 75
$ws.Range("AD29").Value = 8.25
$ws.Range("AE29").Value = 14
$ws.Range("AF29").Value = 10.5
$ws.Range("AH29").Value = 26
$ws.Range("AI29").Value = 37
$ws.Range("AJ29").Value = 700
$ws.Range("G29").Value = 2.35
$ws.Range("H29").Value = 3.15
$ws.Range("I29").Value = 2.9
$ws.Range("L29").Value = 1.35
$ws.Range("M29").Value = 2.72
$ws.Range("N29").Value = 2.02
$ws.Range("O29").Value = 1.62
$ws.Range("P29").Value = 1.45
$ws.Range("Q29").Value = 2.37
$ws.Range("R29").Value = 1.8
$ws.Range("S29").Value = 1.8
$ws.Range("T29").Value = 7.2
$ws.Range("V29").Value = 9.25
$ws.Range("W29").Value = 24
$ws.Range("X29").Value = 21
$ws.Range("Y29").Value = 32
$ws.Range("Z29").Value = 8.5
# Row 38
$ws.Range("AH38").Value = 23
# Row 41
$ws.Range("AA41").Value = 7.2
$ws.Range("AB41").Value = 17.5
$ws.Range("AC41").Value = 90
$ws.Range("AD41").Value = 14.5
$ws.Range("AE41").Value = 35
$ws.Range("AI41").Value = 65
$ws.Range("AJ41").Value = 800
$ws.Range("G41").Value = 1.55
$ws.Range("H41").Value = 3.65
$ws.Range("I41").Value = 5.8
$ws.Range("L41").Value = 1.29
$ws.Range("M41").Value = 3
$ws.Range("N41").Value = 1.85
$ws.Range("O41").Value = 1.75
$ws.Range("P41").Value = 1.4
$ws.Range("Q41").Value = 2.52
$ws.Range("R41").Value = 1.88
$ws.Range("S41").Value = 1.72
$ws.Range("T41").Value = 6.2
$ws.Range("U41").Value = 6.9
$ws.Range("V41").Value = 8
$ws.Range("X41").Value = 13
$ws.Range("Y41").Value = 28
$ws.Range("Z41").Value = 9.5
# Row 42
$ws.Range("AE42").Value = 15.5
$ws.Range("AF42").Value = 11.75
$ws.Range("AH42").Value = 32
$ws.Range("G42").Value = 2.1
$ws.Range("I42").Value = 3.2
$ws.Range("L42").Value = 1.36
$ws.Range("N42").Value = 2.05
$ws.Range("Q42").Value = 2.47
$ws.Range("U42").Value = 9.25
$ws.Range("V42").Value = 9
$ws.Range("W42").Value = 19
$ws.Range("X42").Value = 18.5
$ws.Range("Z42").Value = 8.5
